# Update cryptocurrency price/volume data per the Apr 13 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay a TEXT cell even when it looks numeric
# (e.g. "0.998"), without leaving the cells style pointing at a new index -
# NumberFormat is set to Text just long enough for the assignment to stick as a
# string, then the cell style is reset back to Normal (style index 0).
function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '67.550.27'
$ws.Range("E2").Value = '  -2.60%  '
$ws.Range("D3").Value = '3.264.37'
$ws.Range("E3").Value = '  -5.21%  '
Set-TextValue "D4" '0.998'
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue "D5" '592.29'
$ws.Range("E5").Value = '  -2.81%  '
Set-TextValue "D6" '149.70'
$ws.Range("E6").Value = '  -10.42%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.256.91'
$ws.Range("E8").Value = '  -5.19%  '
Set-TextValue "D9" '0.543'
$ws.Range("E9").Value = '  -8.51%  '
Set-TextValue "D10" '0.171'
$ws.Range("E10").Value = '  -11.40%  '
Set-TextValue "D11" '6.73'
$ws.Range("E11").Value = '  -4.17%  '
Set-TextValue "D12" '0.506'
$ws.Range("E12").Value = '  -10.21%  '
Set-TextValue "D13" '38.51'
$ws.Range("E13").Value = '  -12.95%  '
Set-TextValue "D14" '0.0000247'
$ws.Range("E14").Value = '  -8.57%  '
$ws.Range("D15").Value = '3.778.18'
$ws.Range("E15").Value = '  -5.47%  '
$ws.Range("D16").Value = '67.524.71'
$ws.Range("E16").Value = '  -2.69%  '
$ws.Range("D17").Value = '3.259.02'
$ws.Range("E17").Value = '  -5.57%  '
$ws.Range("E18").Value = '  -5.58%  '
Set-TextValue "D19" '528.94'
$ws.Range("E19").Value = '  -9.02%  '
Set-TextValue "D20" '7.13'
$ws.Range("E20").Value = '  -12.37%  '
Set-TextValue "D21" '14.99'
$ws.Range("E21").Value = '  -12.53%  '
Set-TextValue "D22" '0.756'
$ws.Range("E22").Value = '  -10.63%  '
Set-TextValue "D23" '7.85'
$ws.Range("E23").Value = '  -11.88%  '
Set-TextValue "D24" '85.57'
$ws.Range("E24").Value = '  -11.26%  '
Set-TextValue "D25" '13.51'
$ws.Range("E25").Value = '  -10.71%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -11.16%  '
Set-TextValue "D28" '2.15'
$ws.Range("E28").Value = '  -11.77%  '
Set-TextValue "D29" '8.03'
$ws.Range("E29").Value = '  -7.21%  '
Set-TextValue "D30" '29.02'
$ws.Range("E30").Value = '  -11.44%  '
$ws.Range("E31").Value = '  -3.02%  '
Set-TextValue "D32" '2.67'
$ws.Range("E32").Value = '  -4.52%  '
Set-TextValue "D33" '6.62'
$ws.Range("E33").Value = '  -15.39%  '
Set-TextValue "D34" '5.71'
$ws.Range("E34").Value = '  -12.88%  '
$ws.Range("E35").Value = '  -0.04%  '
Set-TextValue "D36" '512.99'
$ws.Range("E36").Value = '  -11.44%  '
Set-TextValue "D37" '0.0445'
$ws.Range("E37").Value = '  -6.60%  '
Set-TextValue "D38" '53.65'
$ws.Range("E38").Value = '  -4.57%  '
Set-TextValue "D39" '0.0854'
$ws.Range("E39").Value = '  -10.52%  '
Set-TextValue "D40" '8.94'
$ws.Range("E40").Value = '  -15.05%  '
$ws.Range("E41").Value = '  -10.90%  '
Set-TextValue "D42" '2.78'
$ws.Range("E42").Value = '  -11.57%  '
$ws.Range("D43").Value = '2.939.57'
$ws.Range("E43").Value = '  -9.16%  '
Set-TextValue "D44" '0.266'
$ws.Range("E44").Value = '  -9.85%  '
$ws.Range("D45").Value = '0.0₃0588'
$ws.Range("E45").Value = '  -14.29%  '
Set-TextValue "D46" '2.21'
$ws.Range("E46").Value = '  -8.16%  '
Set-TextValue "D47" '26.63'
$ws.Range("E47").Value = '  -14.53%  '
$ws.Range("E48").Value = '  -0.05%  '
Set-TextValue "D49" '2.32'
$ws.Range("E49").Value = '  -16.44%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D50" '0.113'
$ws.Range("E50").Value = '  -10.19%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D51" '123.83'
$ws.Range("E51").Value = '  -7.83%  '
